$wb = $excel.ActiveWorkbook

$ws_ALC = $wb.Worksheets.Item("ALC")
$ws_ARM = $wb.Worksheets.Item("ARM")
$ws_BSM = $wb.Worksheets.Item("BSM")
$ws_CRP = $wb.Worksheets.Item("CRP")
$ws_CUL = $wb.Worksheets.Item("CUL")
$ws_GSM = $wb.Worksheets.Item("GSM")
$ws_LTW = $wb.Worksheets.Item("LTW")
$ws_WVR = $wb.Worksheets.Item("WVR")

# --- ALC ---
$ws_ALC.Range("H132").Value = 1362878.6
$ws_ALC.Range("I132").Value = 1739.4546
$ws_ALC.Range("K132").Value = 5218.3638
$ws_ALC.Range("M132").Value = -2688.3638
$ws_ALC.Range("H137").Value = 2704475.2
$ws_ALC.Range("I137").Value = 3227333.5
$ws_ALC.Range("J137").Value = 3041.6667
$ws_ALC.Range("K137").Value = 9682000.5
$ws_ALC.Range("L137").Value = 9125.000100000001
$ws_ALC.Range("M137").Value = -9679450.5
$ws_ALC.Range("N137").Value = -14225.0001

# --- ARM ---
$ws_ARM.Range("H32").Value = 13687.969
$ws_ARM.Range("I32").Value = 14253.458
$ws_ARM.Range("K32").Value = 14253.458
$ws_ARM.Range("M32").Value = -13966.458
$ws_ARM.Range("H61").Value = 71430320
$ws_ARM.Range("I61").Value = 83334536
$ws_ARM.Range("J61").Value = 4990
$ws_ARM.Range("K61").Value = 83334536
$ws_ARM.Range("L61").Value = 4990
$ws_ARM.Range("M61").Value = -83334324
$ws_ARM.Range("N61").Value = -5414
$ws_ARM.Range("H63").Value = 3419.2222
$ws_ARM.Range("I63").Value = 2994.8
$ws_ARM.Range("J63").Value = 3949.75
$ws_ARM.Range("K63").Value = 2994.8
$ws_ARM.Range("L63").Value = 3949.75
$ws_ARM.Range("M63").Value = -2308.8
$ws_ARM.Range("N63").Value = -5321.75
$ws_ARM.Range("H66").Value = 3419.2222
$ws_ARM.Range("I66").Value = 2994.8
$ws_ARM.Range("J66").Value = 3949.75
$ws_ARM.Range("K66").Value = 14974
$ws_ARM.Range("L66").Value = 19748.75
$ws_ARM.Range("M66").Value = -11542
$ws_ARM.Range("N66").Value = -26612.75
$ws_ARM.Range("H74").Value = 4201487.5
$ws_ARM.Range("I74").Value = 4922708
$ws_ARM.Range("J74").Value = 114572.22
$ws_ARM.Range("K74").Value = 4922708
$ws_ARM.Range("L74").Value = 114572.22
$ws_ARM.Range("M74").Value = -4921834
$ws_ARM.Range("N74").Value = -116320.22
$ws_ARM.Range("H77").Value = 4201487.5
$ws_ARM.Range("I77").Value = 4922708
$ws_ARM.Range("J77").Value = 114572.22
$ws_ARM.Range("K77").Value = 24613540
$ws_ARM.Range("L77").Value = 572861.1
$ws_ARM.Range("M77").Value = -24609172
$ws_ARM.Range("N77").Value = -581597.1
$ws_ARM.Range("H97").Value = 2718100.5
$ws_ARM.Range("I97").Value = 3125607.5
$ws_ARM.Range("J97").Value = 1387
$ws_ARM.Range("K97").Value = 3125607.5
$ws_ARM.Range("L97").Value = 1387
$ws_ARM.Range("M97").Value = -3125111.5
$ws_ARM.Range("N97").Value = -2379
$ws_ARM.Range("H102").Value = 15874760
$ws_ARM.Range("I102").Value = 20409364
$ws_ARM.Range("J102").Value = 3650
$ws_ARM.Range("K102").Value = 20409364
$ws_ARM.Range("L102").Value = 3650
$ws_ARM.Range("M102").Value = -20407742
$ws_ARM.Range("N102").Value = -6894
$ws_ARM.Range("H136").Value = 71430320
$ws_ARM.Range("I136").Value = 83334536
$ws_ARM.Range("J136").Value = 4990
$ws_ARM.Range("K136").Value = 250003608
$ws_ARM.Range("L136").Value = 14970
$ws_ARM.Range("M136").Value = -250001058
$ws_ARM.Range("N136").Value = -20070

# --- BSM ---
$ws_BSM.Range("H134").Value = 5798.6113
$ws_BSM.Range("I134").Value = 5753.222
$ws_BSM.Range("J134").Value = 5844
$ws_BSM.Range("K134").Value = 17259.666
$ws_BSM.Range("L134").Value = 17532
$ws_BSM.Range("M134").Value = -14724.666
$ws_BSM.Range("N134").Value = -22602

# --- CRP ---
$ws_CRP.Range("H6").Value = 10111367
$ws_CRP.Range("I6").Value = 11375163
$ws_CRP.Range("J6").Value = 1000
$ws_CRP.Range("K6").Value = 11375163
$ws_CRP.Range("L6").Value = 1000
$ws_CRP.Range("M6").Value = -11375050
$ws_CRP.Range("N6").Value = -1226
$ws_CRP.Range("H12").Value = 991.25
$ws_CRP.Range("I12").Value = 991.25
$ws_CRP.Range("J12").Value = 0
$ws_CRP.Range("K12").Value = 991.25
$ws_CRP.Range("L12").Value = 0
$ws_CRP.Range("M12").ClearContents()
$ws_CRP.Range("N12").Value = -821.25
$ws_CRP.Range("H16").Value = 1151.625
$ws_CRP.Range("I16").Value = 578.25
$ws_CRP.Range("J16").Value = 1725
$ws_CRP.Range("K16").Value = 578.25
$ws_CRP.Range("L16").Value = 1725
$ws_CRP.Range("M16").Value = -291.25
$ws_CRP.Range("N16").Value = -2299
$ws_CRP.Range("H58").Value = 31251892
$ws_CRP.Range("I58").Value = 33335008
$ws_CRP.Range("J58").Value = 5150.5
$ws_CRP.Range("K58").Value = 33335008
$ws_CRP.Range("L58").Value = 5150.5
$ws_CRP.Range("M58").Value = -33334805
$ws_CRP.Range("N58").Value = -5556.5
$ws_CRP.Range("H99").Value = 2514.0303
$ws_CRP.Range("I99").Value = 2316.125
$ws_CRP.Range("J99").Value = 2700.2942
$ws_CRP.Range("K99").Value = 2316.125
$ws_CRP.Range("L99").Value = 2700.2942
$ws_CRP.Range("M99").Value = -818.125
$ws_CRP.Range("N99").Value = -5696.2942
$ws_CRP.Range("H105").Value = 761.375
$ws_CRP.Range("I105").Value = 761.375
$ws_CRP.Range("J105").Value = 0
$ws_CRP.Range("K105").Value = 761.375
$ws_CRP.Range("L105").Value = 0
$ws_CRP.Range("M105").ClearContents()
$ws_CRP.Range("N105").Value = 985.625
$ws_CRP.Range("H113").Value = 1151.625
$ws_CRP.Range("I113").Value = 578.25
$ws_CRP.Range("J113").Value = 1725
$ws_CRP.Range("K113").Value = 578.25
$ws_CRP.Range("L113").Value = 1725
$ws_CRP.Range("M113").Value = 1591.75
$ws_CRP.Range("N113").Value = -6065
$ws_CRP.Range("H126").Value = 2514.0303
$ws_CRP.Range("I126").Value = 2316.125
$ws_CRP.Range("J126").Value = 2700.2942
$ws_CRP.Range("K126").Value = 6948.375
$ws_CRP.Range("L126").Value = 8100.882599999999
$ws_CRP.Range("M126").Value = -4478.375
$ws_CRP.Range("N126").Value = -13040.8826
$ws_CRP.Range("H132").Value = 36056.734
$ws_CRP.Range("I132").Value = 2485.9092
$ws_CRP.Range("J132").Value = 128376.5
$ws_CRP.Range("K132").Value = 7457.7276
$ws_CRP.Range("L132").Value = 385129.5
$ws_CRP.Range("M132").Value = -4927.7276
$ws_CRP.Range("N132").Value = -390189.5
$ws_CRP.Range("H133").Value = 25000
$ws_CRP.Range("J133").Value = 25000
$ws_CRP.Range("L133").Value = 25000
$ws_CRP.Range("N133").Value = -30060
$ws_CRP.Range("H134").Value = 54478.906
$ws_CRP.Range("I134").Value = 3420.1765
$ws_CRP.Range("J134").Value = 271478.5
$ws_CRP.Range("K134").Value = 10260.5295
$ws_CRP.Range("L134").Value = 814435.5
$ws_CRP.Range("M134").Value = -7725.529500000001
$ws_CRP.Range("N134").Value = -819505.5
$ws_CRP.Range("H136").Value = 31251892
$ws_CRP.Range("I136").Value = 33335008
$ws_CRP.Range("J136").Value = 5150.5
$ws_CRP.Range("K136").Value = 100005024
$ws_CRP.Range("L136").Value = 15451.5
$ws_CRP.Range("M136").Value = -100002474
$ws_CRP.Range("N136").Value = -20551.5

# --- CUL ---
$ws_CUL.Range("H131").Value = 973.4691
$ws_CUL.Range("J131").Value = 1042.1549
$ws_CUL.Range("L131").Value = 3126.4647
$ws_CUL.Range("N131").Value = -13206.4647

# --- GSM ---
$ws_GSM.Range("H102").Value = 1266.6666
$ws_GSM.Range("I102").Value = 0
$ws_GSM.Range("J102").Value = 1266.6666
$ws_GSM.Range("K102").Value = 0
$ws_GSM.Range("L102").ClearContents()
$ws_GSM.Range("M102").Value = 1266.6666
$ws_GSM.Range("N102").Value = -4510.6666
$ws_GSM.Range("H126").Value = 1926.1666
$ws_GSM.Range("I126").Value = 1762.5
$ws_GSM.Range("J126").Value = 2253.5
$ws_GSM.Range("K126").Value = 5287.5
$ws_GSM.Range("L126").Value = 6760.5
$ws_GSM.Range("M126").Value = -2817.5
$ws_GSM.Range("N126").Value = -11700.5
$ws_GSM.Range("H132").Value = 97338.09
$ws_GSM.Range("I132").Value = 78523.234
$ws_GSM.Range("J132").Value = 127912.25
$ws_GSM.Range("K132").Value = 235569.702
$ws_GSM.Range("L132").Value = 383736.75
$ws_GSM.Range("M132").Value = -233039.702
$ws_GSM.Range("N132").Value = -388796.75

# --- LTW ---
$ws_LTW.Range("H7").Value = 2421.0527
$ws_LTW.Range("I7").Value = 2376.4707
$ws_LTW.Range("K7").Value = 2376.4707
$ws_LTW.Range("M7").Value = -2264.4707
$ws_LTW.Range("H22").Value = 1617
$ws_LTW.Range("I22").Value = 0
$ws_LTW.Range("J22").Value = 1617
$ws_LTW.Range("K22").Value = 0
$ws_LTW.Range("L22").ClearContents()
$ws_LTW.Range("M22").Value = 1617
$ws_LTW.Range("N22").Value = -2207
$ws_LTW.Range("H27").Value = 1617
$ws_LTW.Range("I27").Value = 0
$ws_LTW.Range("J27").Value = 1617
$ws_LTW.Range("K27").Value = 0
$ws_LTW.Range("L27").ClearContents()
$ws_LTW.Range("M27").Value = 1617
$ws_LTW.Range("N27").Value = -1831
$ws_LTW.Range("H40").Value = 3080
$ws_LTW.Range("I40").Value = 3080
$ws_LTW.Range("K40").Value = 3080
$ws_LTW.Range("M40").Value = -2944
$ws_LTW.Range("H46").Value = 600.3333
$ws_LTW.Range("I46").Value = 400.5
$ws_LTW.Range("K46").Value = 400.5
$ws_LTW.Range("M46").Value = -212.5
$ws_LTW.Range("H126").Value = 2421.0527
$ws_LTW.Range("I126").Value = 2376.4707
$ws_LTW.Range("K126").Value = 7129.4121
$ws_LTW.Range("M126").Value = -4659.4121
$ws_LTW.Range("H132").Value = 42051.95
$ws_LTW.Range("I132").Value = 19499.172
$ws_LTW.Range("K132").Value = 58497.516
$ws_LTW.Range("M132").Value = -55967.516
$ws_LTW.Range("H136").Value = 90590.69500000001
$ws_LTW.Range("I136").Value = 56865.89
$ws_LTW.Range("J136").Value = 212000
$ws_LTW.Range("K136").Value = 170597.67
$ws_LTW.Range("L136").Value = 636000
$ws_LTW.Range("M136").Value = -168047.67
$ws_LTW.Range("N136").Value = -641100

# --- WVR ---
$ws_WVR.Range("H32").Value = 3000
$ws_WVR.Range("I32").Value = 3000
$ws_WVR.Range("K32").Value = 3000
$ws_WVR.Range("M32").Value = -2683
$ws_WVR.Range("H107").Value = 283
$ws_WVR.Range("I107").Value = 250
$ws_WVR.Range("J107").Value = 316
$ws_WVR.Range("K107").Value = 750
$ws_WVR.Range("L107").Value = 948
$ws_WVR.Range("M107").Value = 1170
$ws_WVR.Range("N107").Value = -4788
$ws_WVR.Range("H123").Value = 45189
$ws_WVR.Range("I123").Value = 42000
$ws_WVR.Range("J123").Value = 49441
$ws_WVR.Range("K123").Value = 42000
$ws_WVR.Range("L123").Value = 49441
$ws_WVR.Range("M123").Value = -37100
$ws_WVR.Range("N123").Value = -59241
$ws_WVR.Range("H126").Value = 928.6
$ws_WVR.Range("I126").Value = 776.6875
$ws_WVR.Range("J126").Value = 1198.6666
$ws_WVR.Range("K126").Value = 2330.0625
$ws_WVR.Range("L126").Value = 3595.9998
$ws_WVR.Range("M126").Value = 139.9375
$ws_WVR.Range("N126").Value = -8535.9998
$ws_WVR.Range("H132").Value = 73072.57000000001
$ws_WVR.Range("I132").Value = 51201.7
$ws_WVR.Range("J132").Value = 127749.75
$ws_WVR.Range("K132").Value = 153605.1
$ws_WVR.Range("L132").Value = 383249.25
$ws_WVR.Range("M132").Value = -151075.1
$ws_WVR.Range("N132").Value = -388309.25
